# Update the task_point data table (fixed dqn & ddqn n_step_update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    1  = @(0.22,   0.473,  0.021,   0, -14,   9)
    2  = @(-0.031, -0.162, 0.604, -178,  52, 142)
    3  = @(0.325,  0.361,  0.32,   -38, -63, 126)
    4  = @(-0.03,  -0.182, 0.107, -156,  58, 102)
    5  = @(0.114,  0.183,  0.459,  -40,  14,   6)
    6  = @(-0.169, -0.03,  0.826,  -89,  40, 103)
    7  = @(-0.121, 0.121,  0.436,   92,  64, -138)
    8  = @(-0.035, 0.173, -0.129,  -92,  50,  -4)
    9  = @(0.018,  -0.172, 0.875,   74, -62,   2)
    10 = @(-0.091, 0.15,   0.011,   53,  18, -163)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
